# Update test table for SkullIsland theme
$wb = $excel.ActiveWorkbook
$normal = $wb.Worksheets.Item(1)

# --- Cosmetic tweaks to the existing "Normal" sheet -----------------------
# Column widths were re-measured (~font metric change); get as close as the
# COM width model allows.
$normal.Range("A1").ColumnWidth = 17.93
$normal.Range("G1:H1").ColumnWidth = 17.5
$normal.Range("I1").ColumnWidth = 40.5

# Row 11 / 19 separator rows grew from 15 to 16 points.
$normal.Rows("11").RowHeight = 16
$normal.Rows("19").RowHeight = 16

# --- Add the new "SkullIsland" worksheet right after "Normal" -------------
$ws = $wb.Worksheets.Add($null, $normal)
$ws.Name = "SkullIsland"

# Header row
$ws.Range("A1").Value = "Skull"
$ws.Range("B1").Value = "Monkey"
$ws.Range("C1").Value = "Parrot"
$ws.Range("D1").Value = "Sword"
$ws.Range("E1").Value = "Coin"
$ws.Range("F1").Value = "Diamond"
$ws.Range("G1").Value = "expect"
$ws.Range("H1").Value = "skull from card"
$ws.Range("I1").Value = "sum check"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 3
$ws.Range("F2").Value = 4
$ws.Range("G2").Formula = "=A2*(-100)"
$ws.Range("I2").Formula = "=SUM(A2:F2)-H2"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("E3").Value = 3

# Rows 4-11 inputs
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("D4").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("H4").Value = 2

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2

$ws.Range("A6").Value = 5
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = 2
$ws.Range("H6").Value = 1

$ws.Range("A7").Value = 6
$ws.Range("E7").Value = 3
$ws.Range("H7").Value = 1

$ws.Range("A8").Value = 7
$ws.Range("E8").Value = 1

$ws.Range("A9").Value = 8

$ws.Range("A10").Value = 9
$ws.Range("H10").Value = 1

$ws.Range("A11").Value = 10
$ws.Range("H11").Value = 2

# Fill formulas down as shared formulas across G3:G11 and I3:I11
$ws.Range("G3:G11").Formula = "=A3*(-100)"
$ws.Range("I3:I11").Formula = "=SUM(A3:F3)-H3"

# Column widths for the new sheet (best match under the COM width model)
$ws.Range("G1:H1").ColumnWidth = 12.5

# Selection / scroll state matching the authored workbook
[void]$ws.Range("F10").Select()
